# Apply the "added the doif rules" edit to the syntactic table workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Append the action-rule label to the DOIF production
#    ( EXPR ) { S } MAYBE_NOT  ->  ( EXPR ) { S } MAYBE_NOT | do_if_rule
$ws.Range("X31").Value = "( EXPR ) { S } MAYBE_NOT | do_if_rule"

# 2. Append the action-rule label to the MAYBE_NOT production
#    maybe_not { S }  ->  maybe_not { S } | maybe_not_rule
$ws.Range("H32").Value = "maybe_not { S } | maybe_not_rule"

# 3. Append the generic "demand_action_rule" action label to every empty
#    (epsilon) production cell belonging to the MAYBE_NOT row (row 32).
$epsilonCells = @("B32","C32","D32","E32","F32","G32","U32","V32","W32","X32","Y32","AB32","AD32")
foreach ($addr in $epsilonCells) {
    $ws.Range($addr).Value = "ϵ | demand_action_rule"
}

# 4. Move the cursor/selection to the cell that was last edited, matching
#    the author's final view in the workbook.
$ws.Range("H32").Select()
